# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a new "Label" column (H) and refreshes the refit prediction /
# error / cross-entropy values in columns D, E and F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Label" header in column H, styled like the other header cells ---
$ws.Range("H1").Value = "Label"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4160
$ws.Range("H1").Borders.LineStyle = 1

# --- Updated refit results for the first (100-iteration) block, rows 2-11 ---
$ws.Range("D2").Value = 0.5234023788613958
$ws.Range("E2").Value = 0.5234023788613958

$ws.Range("D3").Value = 0.3976152680353147
$ws.Range("E3").Value = 0.3976152680353147

$ws.Range("D4").Value = 0.6142760886641699
$ws.Range("E4").Value = 0.6142760886641699

$ws.Range("D5").Value = [double]"9.722172785909883E-18"
$ws.Range("E5").Value = [double]"9.722172785909883E-18"

$ws.Range("D6").Value = 0.7511958070468765
$ws.Range("E6").Value = 0.7511958070468765

$ws.Range("D7").Value = 0.5013276025688661
$ws.Range("E7").Value = 0.4986723974311339

$ws.Range("D8").Value = 0.4941233713025641
$ws.Range("E8").Value = 0.5058766286974359

$ws.Range("D9").Value = 0.4955963855707042
$ws.Range("E9").Value = 0.5044036144292958

$ws.Range("D10").Value = 0.4525796950830527
$ws.Range("E10").Value = 0.5474203049169473

$ws.Range("D11").Value = [double]"1.888175553782187E-05"
$ws.Range("E11").Value = 0.9999811182444622
$ws.Range("F11").Value = 1.735923051834106

# --- New "Label" values (column H) for every data row, both blocks ---
$labels = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 1
}

foreach ($row in $labels.Keys) {
    $ws.Cells.Item($row, 8).Value = $labels[$row]
}
